# Update handback status report timestamps and status text to reflect
# newly generated report (commit: "Generate Report for Handback")

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" column (G), rows 2 and 4
# share the same value (also shared with de-de's Correspond Handoff
# Datetime below) and both move forward together.
$wsOverview.Range("G2").Value = "2016-09-06 22:20:36"
$wsOverview.Range("G4").Value = "2016-09-06 22:20:36"

# zh-cn sheet: Status column (E) for rows 2 and 4 changes from "ht" to "mt";
# Correspond Handoff Datetime (H) and Correspond Handback DateTime (K)
# timestamps advance for the same rows.
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E4").Value = "mt"

$wsZhCn.Range("H2").Value = "2016-09-06 22:20:30"
$wsZhCn.Range("H4").Value = "2016-09-06 22:20:30"

$wsZhCn.Range("K2").Value = "2016-09-06 22:20:48"
$wsZhCn.Range("K4").Value = "2016-09-06 22:20:48"

# de-de sheet: Correspond Handoff Datetime (H) for rows 2 and 4 advances
# (this is the same underlying value as Overview!G2/G4), and Correspond
# Handback DateTime (K) for rows 2 and 4 advances separately.
$wsDeDe.Range("H2").Value = "2016-09-06 22:20:36"
$wsDeDe.Range("H4").Value = "2016-09-06 22:20:36"

$wsDeDe.Range("K2").Value = "2016-09-06 22:20:56"
$wsDeDe.Range("K4").Value = "2016-09-06 22:20:56"
